# Add a "Type" column to the unregistered-students sheet.
#
# Current layout:  A=Name  B=Code  C=Date(+TODAY() formula)
# New layout:       A=Name  B=Code  C=Type(new)  D=Date(+TODAY() formula)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new column: insert a blank column D. Because column D
# was empty, LibreOffice/Excel-style "insert" only creates a new blank
# column there (picking up matching number formats from column C) - the
# existing Date/TODAY() data in column C is left untouched, ready to be
# relocated below.
$ws.Columns("D:D").Insert()

# Move the "Date" header and the TODAY() formulas from column C into the
# freshly inserted column D.
$ws.Range("D1").Value = "Date"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=TODAY()"
}

# Fill column C with the new "Type" header and values. The first four
# students (rows 2-5) are "پارالێل" and the rest (rows 6-9) are "بەیانیان".
$ws.Range("C1").Value = "Type"
$ws.Range("C2:C5").Value = "پارالێل"
$ws.Range("C6:C9").Value = "بەیانیان"

# Keep formatting consistent: apply the header's style across the data
# rows (A2:C9) so the new Type column matches the look of the rest of the
# table.
$headerStyle = $ws.Range("A1").Style
$ws.Range("A2:C9").Style = $headerStyle

# Match the author's final selection/cursor position.
$ws.Range("D11").Select() | Out-Null
